{"js": "// Add user data (Name, NRP, Departement) to the download menu letter.\n// Targets the three empty label lines and appends the corresponding value\n// right after the existing colon, matching the original text layout.\n\nconst updates = [\n  { search: \"         Name            :\", append: \" Patrick\" },\n  { search: \"         NRP              :\", append: \" 05111840000098\" },\n  { search: \"         Departement :\", append: \" Informatika\" },\n];\n\nfor (const { search, append } of updates) {\n  const results = context.document.body.search(search, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text: \"${search}\"`);\n  }\n\n  // Use the first (and expected only) match; insert the value right after it.\n  const range = results.items[0];\n  range.insertText(append, Word.InsertLocation.end);\n}\n\nawait context.sync();\n", "ps1": "# Add user data (Name, NRP, Departement) to the download menu letter.\n$d = $word.ActiveDocument\n\n$updates = @(\n    @{ Find = \"         Name            :\"; Append = \" Patrick\" },\n    @{ Find = \"         NRP              :\"; Append = \" 05111840000098\" },\n    @{ Find = \"         Departement :\"; Append = \" Informatika\" }\n)\n\nforeach ($u in $updates) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.Text = $u.Find\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute()\n\n    if ($find.Found) {\n        $range.Collapse(0)  # wdCollapseEnd\n        $range.InsertAfter($u.Append)\n    }\n}\n"}
